# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-06 20:16:40
#
# The "Recorded By" column (G) on the active sheet lists the users who
# recorded/edited attendance as a comma separated string (e.g.
# "System, dnasr281@gmail.com"). For a specific set of rows the two
# names were reordered (the second name moved to the front / the first
# name moved to the back). This script reproduces that reordering by
# swapping the two comma separated parts for each of the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-based worksheet rows) whose column G value needs its two
# comma-separated entries swapped, taken from the diff.
$rows = @(3,4,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,30,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,56,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7 ("Recorded By")
    $value = $cell.Value2

    if ($value -ne $null) {
        $parts = $value -split ",\s*"

        if ($parts.Count -eq 2) {
            $newValue = ($parts[1].Trim()) + ", " + ($parts[0].Trim())
            $cell.Value2 = $newValue
        }
    }
}
